$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header "Save" in H1, matching the style of the other header cells (e.g. G1)
# by copying G1's format (bold, centered, bordered) onto H1, then setting its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add value 0 in H2 (plain numeric cell, no special style - like F2/G2)
$ws.Range("H2").Value = 0
